$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 298 (id=296)
$r = 298
$ws.Range("A2:AB2").Copy($ws.Range("A" + $r + ":AB" + $r))
$ws.Cells.Item($r, 1).Value = 296
$ws.Cells.Item($r, 2).Value = 7090293
$ws.Cells.Item($r, 3).Value = "Poland Ekstraklasa"
$ws.Cells.Item($r, 4).Value = 45437.52083333334
$ws.Cells.Item($r, 5).Value = "Radomiak Radom"
$ws.Cells.Item($r, 6).Value = "Widzew Lodz"
$ws.Cells.Item($r, 7).Value = 1
$ws.Cells.Item($r, 8).Value = 3
$ws.Cells.Item($r, 9).Value = "A"
$ws.Cells.Item($r, 10).Value = 2.2
$ws.Cells.Item($r, 11).Value = 3.1
$ws.Cells.Item($r, 12).Value = 3.1
$ws.Cells.Item($r, 13).Value = 2.15
$ws.Cells.Item($r, 14).Value = 3.2
$ws.Cells.Item($r, 15).Value = 3.1
$ws.Cells.Item($r, 16).Value = -0.25
$ws.Cells.Item($r, 17).Value = 1.925
$ws.Cells.Item($r, 18).Value = 1.925
$ws.Cells.Item($r, 19).Value = 2.75
$ws.Cells.Item($r, 20).Value = 1.9
$ws.Cells.Item($r, 21).Value = 1.95
$ws.Cells.Item($r, 22).Value = -1
$ws.Cells.Item($r, 23).Value = -1
$ws.Cells.Item($r, 24).Value = 2.1
$ws.Cells.Item($r, 25).Value = -1
$ws.Cells.Item($r, 26).Value = 0.925
$ws.Cells.Item($r, 27).Value = 0.8999999999999999
$ws.Cells.Item($r, 28).Value = -1

# Row 299 (id=297)
$r = 299
$ws.Range("A2:AB2").Copy($ws.Range("A" + $r + ":AB" + $r))
$ws.Cells.Item($r, 1).Value = 297
$ws.Cells.Item($r, 2).Value = 7074364
$ws.Cells.Item($r, 3).Value = "Poland Ekstraklasa"
$ws.Cells.Item($r, 4).Value = 45437.52083333334
$ws.Cells.Item($r, 5).Value = "Rakow Czestochowa"
$ws.Cells.Item($r, 6).Value = "Slask Wroclaw"
$ws.Cells.Item($r, 7).Value = 1
$ws.Cells.Item($r, 8).Value = 2
$ws.Cells.Item($r, 9).Value = "A"
$ws.Cells.Item($r, 10).Value = 2.5
$ws.Cells.Item($r, 11).Value = 3.6
$ws.Cells.Item($r, 12).Value = 2.4
$ws.Cells.Item($r, 13).Value = 2.15
$ws.Cells.Item($r, 14).Value = 3.6
$ws.Cells.Item($r, 15).Value = 2.875
$ws.Cells.Item($r, 16).Value = -0.25
$ws.Cells.Item($r, 17).Value = 1.95
$ws.Cells.Item($r, 18).Value = 1.9
$ws.Cells.Item($r, 19).Value = 2.5
$ws.Cells.Item($r, 20).Value = 1.875
$ws.Cells.Item($r, 21).Value = 1.975
$ws.Cells.Item($r, 22).Value = -1
$ws.Cells.Item($r, 23).Value = -1
$ws.Cells.Item($r, 24).Value = 1.875
$ws.Cells.Item($r, 25).Value = -1
$ws.Cells.Item($r, 26).Value = 0.8999999999999999
$ws.Cells.Item($r, 27).Value = 0.875
$ws.Cells.Item($r, 28).Value = -1

# Row 300 (id=298)
$r = 300
$ws.Range("A2:AB2").Copy($ws.Range("A" + $r + ":AB" + $r))
$ws.Cells.Item($r, 1).Value = 298
$ws.Cells.Item($r, 2).Value = 7041338
$ws.Cells.Item($r, 3).Value = "Poland Ekstraklasa"
$ws.Cells.Item($r, 4).Value = 45437.52083333334
$ws.Cells.Item($r, 5).Value = "Jagiellonia Bialystok"
$ws.Cells.Item($r, 6).Value = "Warta Poznan"
$ws.Cells.Item($r, 7).Value = 3
$ws.Cells.Item($r, 8).Value = 0
$ws.Cells.Item($r, 9).Value = "H"
$ws.Cells.Item($r, 10).Value = 1.444
$ws.Cells.Item($r, 11).Value = 4.75
$ws.Cells.Item($r, 12).Value = 5.25
$ws.Cells.Item($r, 13).Value = 1.4
$ws.Cells.Item($r, 14).Value = 4.75
$ws.Cells.Item($r, 15).Value = 5.75
$ws.Cells.Item($r, 16).Value = -1.25
$ws.Cells.Item($r, 17).Value = 1.9
$ws.Cells.Item($r, 18).Value = 1.95
$ws.Cells.Item($r, 19).Value = 3
$ws.Cells.Item($r, 20).Value = 1.925
$ws.Cells.Item($r, 21).Value = 1.925
$ws.Cells.Item($r, 22).Value = 0.3999999999999999
$ws.Cells.Item($r, 23).Value = -1
$ws.Cells.Item($r, 24).Value = -1
$ws.Cells.Item($r, 25).Value = 0.8999999999999999
$ws.Cells.Item($r, 26).Value = -1
$ws.Cells.Item($r, 27).Value = 0
$ws.Cells.Item($r, 28).Value = 0

# Row 301 (id=299)
$r = 301
$ws.Range("A2:AB2").Copy($ws.Range("A" + $r + ":AB" + $r))
$ws.Cells.Item($r, 1).Value = 299
$ws.Cells.Item($r, 2).Value = 7083187
$ws.Cells.Item($r, 3).Value = "Poland Ekstraklasa"
$ws.Cells.Item($r, 4).Value = 45437.52083333334
$ws.Cells.Item($r, 5).Value = "Lech Poznan"
$ws.Cells.Item($r, 6).Value = "Korona Kielce"
$ws.Cells.Item($r, 7).Value = 1
$ws.Cells.Item($r, 8).Value = 2
$ws.Cells.Item($r, 9).Value = "A"
$ws.Cells.Item($r, 10).Value = 1.8
$ws.Cells.Item($r, 11).Value = 3.8
$ws.Cells.Item($r, 12).Value = 3.6
$ws.Cells.Item($r, 13).Value = 2.1
$ws.Cells.Item($r, 14).Value = 3.7
$ws.Cells.Item($r, 15).Value = 2.9
$ws.Cells.Item($r, 16).Value = -0.25
$ws.Cells.Item($r, 17).Value = 1.9
$ws.Cells.Item($r, 18).Value = 1.95
$ws.Cells.Item($r, 19).Value = 2.75
$ws.Cells.Item($r, 20).Value = 1.925
$ws.Cells.Item($r, 21).Value = 1.925
$ws.Cells.Item($r, 22).Value = -1
$ws.Cells.Item($r, 23).Value = -1
$ws.Cells.Item($r, 24).Value = 1.9
$ws.Cells.Item($r, 25).Value = -1
$ws.Cells.Item($r, 26).Value = 0.95
$ws.Cells.Item($r, 27).Value = 0.4625
$ws.Cells.Item($r, 28).Value = -0.5

# Row 302 (id=300)
$r = 302
$ws.Range("A2:AB2").Copy($ws.Range("A" + $r + ":AB" + $r))
$ws.Cells.Item($r, 1).Value = 300
$ws.Cells.Item($r, 2).Value = 7088350
$ws.Cells.Item($r, 3).Value = "Poland Ekstraklasa"
$ws.Cells.Item($r, 4).Value = 45437.52083333334
$ws.Cells.Item($r, 5).Value = "Puszcza Niepolomice"
$ws.Cells.Item($r, 6).Value = "Piast Gliwice"
$ws.Cells.Item($r, 7).Value = 1
$ws.Cells.Item($r, 8).Value = 0
$ws.Cells.Item($r, 9).Value = "H"
$ws.Cells.Item($r, 10).Value = 3
$ws.Cells.Item($r, 11).Value = 3.1
$ws.Cells.Item($r, 12).Value = 2.3
$ws.Cells.Item($r, 13).Value = 2.7
$ws.Cells.Item($r, 14).Value = 3
$ws.Cells.Item($r, 15).Value = 2.625
$ws.Cells.Item($r, 16).Value = 0
$ws.Cells.Item($r, 17).Value = 1.975
$ws.Cells.Item($r, 18).Value = 1.875
$ws.Cells.Item($r, 19).Value = 2.25
$ws.Cells.Item($r, 20).Value = 2.025
$ws.Cells.Item($r, 21).Value = 1.825
$ws.Cells.Item($r, 22).Value = 1.7
$ws.Cells.Item($r, 23).Value = -1
$ws.Cells.Item($r, 24).Value = -1
$ws.Cells.Item($r, 25).Value = 0.9750000000000001
$ws.Cells.Item($r, 26).Value = -1
$ws.Cells.Item($r, 27).Value = -1
$ws.Cells.Item($r, 28).Value = 0.825

# Row 303 (id=301)
$r = 303
$ws.Range("A2:AB2").Copy($ws.Range("A" + $r + ":AB" + $r))
$ws.Cells.Item($r, 1).Value = 301
$ws.Cells.Item($r, 2).Value = 7083189
$ws.Cells.Item($r, 3).Value = "Poland Ekstraklasa"
$ws.Cells.Item($r, 4).Value = 45437.52083333334
$ws.Cells.Item($r, 5).Value = "Pogon Szczecin"
$ws.Cells.Item($r, 6).Value = "Gornik Zabrze"
$ws.Cells.Item($r, 7).Value = 1
$ws.Cells.Item($r, 8).Value = 0
$ws.Cells.Item($r, 9).Value = "H"
$ws.Cells.Item($r, 10).Value = 1.727
$ws.Cells.Item($r, 11).Value = 4
$ws.Cells.Item($r, 12).Value = 3.75
$ws.Cells.Item($r, 13).Value = 1.55
$ws.Cells.Item($r, 14).Value = 4.333
$ws.Cells.Item($r, 15).Value = 4.5
$ws.Cells.Item($r, 16).Value = -1
$ws.Cells.Item($r, 17).Value = 1.925
$ws.Cells.Item($r, 18).Value = 1.925
$ws.Cells.Item($r, 19).Value = 3.5
$ws.Cells.Item($r, 20).Value = 2.025
$ws.Cells.Item($r, 21).Value = 1.825
$ws.Cells.Item($r, 22).Value = 0.55
$ws.Cells.Item($r, 23).Value = -1
$ws.Cells.Item($r, 24).Value = -1
$ws.Cells.Item($r, 25).Value = 0
$ws.Cells.Item($r, 26).Value = 0
$ws.Cells.Item($r, 27).Value = -1
$ws.Cells.Item($r, 28).Value = 0.825

# Row 304 (id=302)
$r = 304
$ws.Range("A2:AB2").Copy($ws.Range("A" + $r + ":AB" + $r))
$ws.Cells.Item($r, 1).Value = 302
$ws.Cells.Item($r, 2).Value = 7083188
$ws.Cells.Item($r, 3).Value = "Poland Ekstraklasa"
$ws.Cells.Item($r, 4).Value = 45437.52083333334
$ws.Cells.Item($r, 5).Value = "Legia Warsaw"
$ws.Cells.Item($r, 6).Value = "Zaglebie Lubin"
$ws.Cells.Item($r, 7).Value = 2
$ws.Cells.Item($r, 8).Value = 1
$ws.Cells.Item($r, 9).Value = "H"
$ws.Cells.Item($r, 10).Value = 1.5
$ws.Cells.Item($r, 11).Value = 4
$ws.Cells.Item($r, 12).Value = 5.5
$ws.Cells.Item($r, 13).Value = 1.6
$ws.Cells.Item($r, 14).Value = 4.1
$ws.Cells.Item($r, 15).Value = 4.333
$ws.Cells.Item($r, 16).Value = -0.75
$ws.Cells.Item($r, 17).Value = 1.825
$ws.Cells.Item($r, 18).Value = 2.025
$ws.Cells.Item($r, 19).Value = 3
$ws.Cells.Item($r, 20).Value = 1.875
$ws.Cells.Item($r, 21).Value = 1.975
$ws.Cells.Item($r, 22).Value = 0.6000000000000001
$ws.Cells.Item($r, 23).Value = -1
$ws.Cells.Item($r, 24).Value = -1
$ws.Cells.Item($r, 25).Value = 0.4125
$ws.Cells.Item($r, 26).Value = -0.5
$ws.Cells.Item($r, 27).Value = 0
$ws.Cells.Item($r, 28).Value = 0

# Row 305 (id=303)
$r = 305
$ws.Range("A2:AB2").Copy($ws.Range("A" + $r + ":AB" + $r))
$ws.Cells.Item($r, 1).Value = 303
$ws.Cells.Item($r, 2).Value = 7093821
$ws.Cells.Item($r, 3).Value = "Poland Ekstraklasa"
$ws.Cells.Item($r, 4).Value = 45437.52083333334
$ws.Cells.Item($r, 5).Value = "LKS Lodz"
$ws.Cells.Item($r, 6).Value = "Stal Mielec"
$ws.Cells.Item($r, 7).Value = 3
$ws.Cells.Item($r, 8).Value = 2
$ws.Cells.Item($r, 9).Value = "H"
$ws.Cells.Item($r, 10).Value = 2.5
$ws.Cells.Item($r, 11).Value = 3.4
$ws.Cells.Item($r, 12).Value = 2.5
$ws.Cells.Item($r, 13).Value = 2.2
$ws.Cells.Item($r, 14).Value = 3.5
$ws.Cells.Item($r, 15).Value = 2.8
$ws.Cells.Item($r, 16).Value = -0.25
$ws.Cells.Item($r, 17).Value = 2.025
$ws.Cells.Item($r, 18).Value = 1.825
$ws.Cells.Item($r, 19).Value = 3
$ws.Cells.Item($r, 20).Value = 2
$ws.Cells.Item($r, 21).Value = 1.85
$ws.Cells.Item($r, 22).Value = 1.2
$ws.Cells.Item($r, 23).Value = -1
$ws.Cells.Item($r, 24).Value = -1
$ws.Cells.Item($r, 25).Value = 1.025
$ws.Cells.Item($r, 26).Value = -1
$ws.Cells.Item($r, 27).Value = 1
$ws.Cells.Item($r, 28).Value = -1

# Row 306 (id=304)
$r = 306
$ws.Range("A2:AB2").Copy($ws.Range("A" + $r + ":AB" + $r))
$ws.Cells.Item($r, 1).Value = 304
$ws.Cells.Item($r, 2).Value = 7093820
$ws.Cells.Item($r, 3).Value = "Poland Ekstraklasa"
$ws.Cells.Item($r, 4).Value = 45437.52083333334
$ws.Cells.Item($r, 5).Value = "Ruch Chorzow"
$ws.Cells.Item($r, 6).Value = "Cracovia Krakow"
$ws.Cells.Item($r, 7).Value = 2
$ws.Cells.Item($r, 8).Value = 0
$ws.Cells.Item($r, 9).Value = "H"
$ws.Cells.Item($r, 10).Value = 2.5
$ws.Cells.Item($r, 11).Value = 3.4
$ws.Cells.Item($r, 12).Value = 2.5
$ws.Cells.Item($r, 13).Value = 2.6
$ws.Cells.Item($r, 14).Value = 3.6
$ws.Cells.Item($r, 15).Value = 2.3
$ws.Cells.Item($r, 16).Value = 0
$ws.Cells.Item($r, 17).Value = 2.025
$ws.Cells.Item($r, 18).Value = 1.825
$ws.Cells.Item($r, 19).Value = 3
$ws.Cells.Item($r, 20).Value = 2.025
$ws.Cells.Item($r, 21).Value = 1.825
$ws.Cells.Item($r, 22).Value = 1.6
$ws.Cells.Item($r, 23).Value = -1
$ws.Cells.Item($r, 24).Value = -1
$ws.Cells.Item($r, 25).Value = 1.025
$ws.Cells.Item($r, 26).Value = -1
$ws.Cells.Item($r, 27).Value = -1
$ws.Cells.Item($r, 28).Value = 0.825
